$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2087227414330218
$ws.Range("C2").Value = 0.5171339563862928
$ws.Range("J2").Value = 0.01869158878504673
$ws.Range("P2").Value = 0.1339563862928349
$ws.Range("S2").Value = 0.1214953271028037

$ws.Range("B3").Value = 0.005917159763313609
$ws.Range("C3").Value = 0.02958579881656805
$ws.Range("J3").Value = 0.01775147928994083
$ws.Range("P3").Value = 0.7041420118343196
$ws.Range("S3").Value = 0.242603550295858

$ws.Range("J4").Value = 0.06382978723404255
$ws.Range("P4").Value = 0.6595744680851063
$ws.Range("S4").Value = 0.2765957446808511

$ws.Range("B6").Value = 0.04741379310344827
$ws.Range("D6").Value = 0.02155172413793104
$ws.Range("F6").Value = 0.03017241379310345
$ws.Range("J6").Value = 0.2198275862068965
$ws.Range("O6").Value = 0.02155172413793104
$ws.Range("Q6").Value = 0.2025862068965517
$ws.Range("R6").Value = 0.06465517241379311
$ws.Range("S6").Value = 0.3922413793103448

$ws.Range("B7").Value = 0.08947368421052632
$ws.Range("D7").Value = 0.03157894736842105
$ws.Range("F7").Value = 0.06315789473684211
$ws.Range("J7").Value = 0.1105263157894737
$ws.Range("O7").Value = 0.02631578947368421
$ws.Range("Q7").Value = 0.1526315789473684
$ws.Range("R7").Value = 0.08421052631578947
$ws.Range("S7").Value = 0.4421052631578947

$ws.Range("B8").Value = 0.07962962962962963
$ws.Range("D8").Value = 0.007407407407407408
$ws.Range("F8").Value = 0.07777777777777778
$ws.Range("J8").Value = 0.09444444444444444
$ws.Range("O8").Value = 0.01296296296296296
$ws.Range("Q8").Value = 0.1944444444444444
$ws.Range("R8").Value = 0.06111111111111111
$ws.Range("S8").Value = 0.4722222222222222

$ws.Range("B9").Value = 0.120253164556962
$ws.Range("D9").Value = 0.03164556962025317
$ws.Range("F9").Value = 0.05696202531645569
$ws.Range("J9").Value = 0.1075949367088608
$ws.Range("O9").Value = 0.006329113924050633
$ws.Range("Q9").Value = 0.1772151898734177
$ws.Range("R9").Value = 0.06962025316455696
$ws.Range("S9").Value = 0.4303797468354431

$ws.Range("B10").Value = 0.1171702284450995
$ws.Range("D10").Value = 0.02137067059690494
$ws.Range("E10").Value = 0.001473839351510685
$ws.Range("F10").Value = 0.079587324981577
$ws.Range("J10").Value = 0.1061164333087693
$ws.Range("O10").Value = 0.01473839351510685
$ws.Range("Q10").Value = 0.2225497420781135
$ws.Range("R10").Value = 0.06779661016949153
$ws.Range("S10").Value = 0.3691967575534267

$ws.Range("G11").Value = 0.1622516556291391
$ws.Range("J11").Value = 0.08609271523178808
$ws.Range("K11").Value = 0.2086092715231788
$ws.Range("L11").Value = 0.5264900662251656
$ws.Range("S11").Value = 0.01655629139072848

$ws.Range("G12").Value = 0.7901234567901234
$ws.Range("J12").Value = 0.1790123456790123
$ws.Range("K12").Value = 0.01234567901234568
$ws.Range("L12").Value = 0.01234567901234568
$ws.Range("S12").Value = 0.006172839506172839

$ws.Range("G13").Value = 0.631578947368421
$ws.Range("J13").Value = 0.3157894736842105
$ws.Range("S13").Value = 0.05263157894736842

$ws.Range("J14").Value = 1

$ws.Range("F15").Value = 0.01754385964912281
$ws.Range("H15").Value = 0.2324561403508772
$ws.Range("I15").Value = 0.03508771929824561
$ws.Range("J15").Value = 0.3728070175438596
$ws.Range("K15").Value = 0.06140350877192982
$ws.Range("M15").Value = 0.008771929824561403
$ws.Range("O15").Value = 0.06578947368421052
$ws.Range("S15").Value = 0.206140350877193

$ws.Range("F16").Value = 0.03191489361702127
$ws.Range("H16").Value = 0.2287234042553191
$ws.Range("I16").Value = 0.04787234042553191
$ws.Range("J16").Value = 0.4361702127659575
$ws.Range("K16").Value = 0.0797872340425532
$ws.Range("O16").Value = 0.05851063829787234
$ws.Range("S16").Value = 0.1170212765957447

$ws.Range("F17").Value = 0.007889546351084813
$ws.Range("H17").Value = 0.1854043392504931
$ws.Range("I17").Value = 0.08481262327416174
$ws.Range("J17").Value = 0.4694280078895464
$ws.Range("K17").Value = 0.08875739644970414
$ws.Range("M17").Value = 0.01183431952662722
$ws.Range("O17").Value = 0.05719921104536489
$ws.Range("S17").Value = 0.09467455621301775

$ws.Range("F18").Value = 0.01212121212121212
$ws.Range("H18").Value = 0.2121212121212121
$ws.Range("I18").Value = 0.08484848484848485
$ws.Range("J18").Value = 0.4484848484848485
$ws.Range("K18").Value = 0.09090909090909091
$ws.Range("M18").Value = 0.01212121212121212
$ws.Range("O18").Value = 0.06060606060606061
$ws.Range("S18").Value = 0.07878787878787878

$ws.Range("F19").Value = 0.01124437781109445
$ws.Range("H19").Value = 0.2353823088455772
$ws.Range("I19").Value = 0.06371814092953523
$ws.Range("J19").Value = 0.3980509745127436
$ws.Range("K19").Value = 0.1079460269865068
$ws.Range("M19").Value = 0.02173913043478261
$ws.Range("N19").Value = 0.0007496251874062968
$ws.Range("O19").Value = 0.06971514242878561
$ws.Range("S19").Value = 0.09145427286356822

Write-Output "Applied team specific time data updates"
